# repull data, push all data, mean calculation
# Update column F (dSF) values for several rows to reflect repulled data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = 3
    3  = -3
    4  = -1
    6  = 1
    8  = -1
    9  = -1
    11 = -5
    12 = -2
    13 = 4
    15 = -4
    16 = -3
    17 = 4
    18 = 2
    19 = 2
    20 = -7
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 6).Value = $updates[$row]
}
